# Update cryptocurrency price (D) and volume-change (E) columns on Sheet1
# to reflect the refreshed values from the automated data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "59.373.89"
$ws.Range("E2").Value = "  -2.27%  "
$ws.Range("D3").Value = "2.579.29"
$ws.Range("E3").Value = "  -3.07%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.65%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("D9").Value = "2.586.67"
$ws.Range("E9").Value = "  -2.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.51%  "
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("E12").Value = "  +12.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.352"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.00%  "
$ws.Range("D14").Value = "3.035.74"
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("D15").Value = "59.375.44"
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.10%  "
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "2.589.04"
$ws.Range("E18").Value = "  -1.84%  "
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "337.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.470"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").Value = "0.0₃0777"
$ws.Range("E29").Value = "  -3.12%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("E32").Value = "  -2.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "157.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.908"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.858"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "291.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "136.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("E45").Value = "  -1.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.593"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0532"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.21%  "
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").Value = "1.954.37"
$ws.Range("E51").Value = "  -0.33%  "
